# Append one new log row (row 83) to each of the four sheets, mirroring the
# previous day's row (row 82) but stamped with the newest timestamp. This
# also extends each sheet's used range from A1:I82 to A1:I83.

$wb = $excel.ActiveWorkbook

$newRow = 83
$timestamp = 45869.46447916667
$dateFormat = "YYYY-MM-DD HH:MM:SS"

$rowsData = @{
    "MID_LFT_#1" = @("0x01,0x90", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,", "0x01,0x30", "0x07", 400, "5.68631262647113e+23", 304, 7)
    "MID_LFT_#2" = @("0x01,0x7c", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,", "0x01,0x30", "0x19", 380, "5.68432987514711e+23", 304, 25)
    "MID_PLT_#1" = @("0x00,0x6e", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,", "0x00,0x60", "0x15", 110, "5.68631262647113e+23", 96, 15)
    "MID_PLT_#2" = @("0x00,0x82", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,", "0x00,0x76", "0x9", 130, "5.68631262647113e+23", 118, 9)
}

foreach ($ws in $wb.Worksheets) {
    $data = $rowsData[$ws.Name]
    if ($data -eq $null) { continue }

    $ws.Cells.Item($newRow, 1).Value = $timestamp
    $ws.Cells.Item($newRow, 1).NumberFormat = $dateFormat

    $ws.Cells.Item($newRow, 2).Value = $data[0]
    $ws.Cells.Item($newRow, 3).Value = $data[1]
    $ws.Cells.Item($newRow, 4).Value = $data[2]
    $ws.Cells.Item($newRow, 5).Value = $data[3]
    $ws.Cells.Item($newRow, 6).Value = $data[4]
    $ws.Cells.Item($newRow, 7).Value = [double]$data[5]
    $ws.Cells.Item($newRow, 8).Value = $data[6]
    $ws.Cells.Item($newRow, 9).Value = $data[7]
}
